# Auto-generated script to apply the "Add data for 2024-09-04" update
# to the violent-crime-full-year workbook. For each affected worksheet,
# the 2024 (column K) and, in a few cases, 2023 (column J) year-to-date
# totals are updated to reflect the newly added day of data.

$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5434
$ws.Range("K3").Value = 5593
$ws.Range("J4").Value = 1832
$ws.Range("K4").Value = 1166
$ws.Range("K5").Value = 400
$ws.Range("K6").Value = 6213
$ws.Range("J7").Value = 29298
$ws.Range("K7").Value = 18806

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K5").Value = 45
$ws.Range("K6").Value = 136
$ws.Range("K7").Value = 556
$ws.Range("K8").Value = 1254
$ws.Range("K9").Value = 77
$ws.Range("K11").Value = 357
$ws.Range("K15").Value = 188
$ws.Range("K19").Value = 550
$ws.Range("K20").Value = 435
$ws.Range("K24").Value = 56
$ws.Range("K25").Value = 88
$ws.Range("K27").Value = 177
$ws.Range("K29").Value = 1013
$ws.Range("K30").Value = 75
$ws.Range("K31").Value = 202
$ws.Range("K33").Value = 804
$ws.Range("K36").Value = 248
$ws.Range("K37").Value = 632
$ws.Range("K39").Value = 24
$ws.Range("K40").Value = 43
$ws.Range("K42").Value = 702
$ws.Range("K44").Value = 162
$ws.Range("K46").Value = 40
$ws.Range("K47").Value = 129
$ws.Range("K48").Value = 238
$ws.Range("K50").Value = 88
$ws.Range("K51").Value = 235
$ws.Range("K52").Value = 491
$ws.Range("K53").Value = 239
$ws.Range("K54").Value = 365
$ws.Range("K55").Value = 210
$ws.Range("K57").Value = 72
$ws.Range("J63").Value = 113
$ws.Range("K65").Value = 429
$ws.Range("K67").Value = 715
$ws.Range("K76").Value = 260
$ws.Range("K77").Value = 131
$ws.Range("K78").Value = 219
$ws.Range("K79").Value = 474
$ws.Range("K80").Value = 65
$ws.Range("K83").Value = 419
$ws.Range("K84").Value = 143
$ws.Range("K85").Value = 888
$ws.Range("K89").Value = 274
$ws.Range("K90").Value = 171
$ws.Range("K91").Value = 209
$ws.Range("K92").Value = 70
$ws.Range("K95").Value = 319
$ws.Range("K96").Value = 202
$ws.Range("K99").Value = 315
$ws.Range("J101").Value = 29298
$ws.Range("K101").Value = 18806

# West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 202

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 179
$ws.Range("K7").Value = 556

# Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 121
$ws.Range("K7").Value = 357

# Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 77
$ws.Range("K3").Value = 84
$ws.Range("K7").Value = 274

# South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 302
$ws.Range("K6").Value = 217
$ws.Range("K7").Value = 888

# Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 136
$ws.Range("K7").Value = 491

# Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 239

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 344
$ws.Range("K3").Value = 379
$ws.Range("K6").Value = 423
$ws.Range("K7").Value = 1254

# South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 142
$ws.Range("K7").Value = 419

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 220
$ws.Range("K7").Value = 804

# West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 319

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 181
$ws.Range("K3").Value = 208
$ws.Range("K7").Value = 632

# New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 138
$ws.Range("K4").Value = 16
$ws.Range("K6").Value = 156
$ws.Range("K7").Value = 429

# Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 315

# Fuller Park
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 75

# Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 202

# North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 204
$ws.Range("K3").Value = 257
$ws.Range("K6").Value = 198
$ws.Range("K7").Value = 715

# South Deering
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 143

# Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 91
$ws.Range("K6").Value = 195
$ws.Range("K7").Value = 365

# Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 288
$ws.Range("K3").Value = 363
$ws.Range("K4").Value = 49
$ws.Range("K6").Value = 286
$ws.Range("K7").Value = 1013

# Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 238

# Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 164
$ws.Range("K7").Value = 550

# Irving Park
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 39
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 162

# River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 56
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 260

# Ashburn
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 50
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 136

# Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 186
$ws.Range("K3").Value = 217
$ws.Range("K4").Value = 28
$ws.Range("K6").Value = 264
$ws.Range("K7").Value = 702

# Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 65
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 219

# Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 64
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 210

# Dunning
$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 56

# Jefferson Park
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 40

# Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 53
$ws.Range("K7").Value = 209

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 474

# Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 138
$ws.Range("K7").Value = 435

# Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 98
$ws.Range("K7").Value = 248

# East Side
$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 88

# Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 129

# Brighton Park
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 188

# Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 88

# Greektown
$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("K3").Value = 1
$ws.Range("K6").Value = 24

# Avalon Park
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K2").Value = 26
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 77

# West Elsdon
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 70

# Armour Square
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 45

# Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K3").Value = 43
$ws.Range("K7").Value = 177

# Washington Heights
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 61
$ws.Range("K7").Value = 171

# Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K4").Value = 25
$ws.Range("K7").Value = 235

# Mckinley Park
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 72

# Riverdale
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 131

# Rush & Division
$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 65

# Hegewisch
$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 43
